$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B values for rows that keep the same label (A) but change B
$ws.Range("B2").Value = 0.4346677117383549
$ws.Range("B3").Value = 0.2378829424012432
$ws.Range("B4").Value = 0.09586447995276519

# Update rows 5-9: both label (A) and value (B) change / reorder
$ws.Range("A5").Value = "VIX_short"
$ws.Range("B5").Value = 0.04825445431697448

$ws.Range("A6").Value = "VIX_long"
$ws.Range("B6").Value = 0.04275628168602999

$ws.Range("A7").Value = "close_short"
$ws.Range("B7").Value = 0.04007826965909733

$ws.Range("A8").Value = "close_long"
$ws.Range("B8").Value = 0.03999035233643641

$ws.Range("A9").Value = "VIX"
$ws.Range("B9").Value = 0.03180555551077596

# Row 10: label stays the same, only value changes
$ws.Range("B10").Value = 0.02869995239832254
